$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new data adds a 3rd sending cluster ("ECs") to the Bmp7->Bmpr2 table.
# Row layout changes from 2 clusters x 4 targets (rows 2-9) to 3 clusters x 4
# targets (rows 2-13): ECs (2-5), FAPs (6-9, shifted down from 2-5), MuSCs
# (10-13, shifted down from 6-9) -- all with refreshed TPM-derived values.
# Append 4 blank rows at the bottom (away from the styled header row) so no
# formatting gets inherited/duplicated, then overwrite every data cell in
# place with its final value.
$ws.Rows("10:13").Insert()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp7"
$ws.Range("C2").Value = "Bmpr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03443933333333333
$ws.Range("H2").Value = 0.103318
$ws.Range("I2").Value = 0.05823261822459219
$ws.Range("J2").Value = 0.0582326182245922
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 46.29121633333333
$ws.Range("N2").Value = 138.873649
$ws.Range("O2").Value = 0.3133663986859022
$ws.Range("P2").Value = 0.3133663986859022
$ws.Range("Q2").Value = 1.594238629709111
$ws.Range("R2").Value = 14.348147667382
$ws.Range("S2").Value = 0.01824814585909149
$ws.Range("T2").Value = 0.01824814585909149

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp7"
$ws.Range("C3").Value = "Bmpr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03443933333333333
$ws.Range("H3").Value = 0.103318
$ws.Range("I3").Value = 0.05823261822459219
$ws.Range("J3").Value = 0.0582326182245922
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("N3").Value = 140.44867
$ws.Range("O3").Value = 0.3169204109998198
$ws.Range("P3").Value = 0.3169204109998198
$ws.Range("Q3").Value = 1.612319520784444
$ws.Range("R3").Value = 14.51087568706
$ws.Range("S3").Value = 0.01845510530133335
$ws.Range("T3").Value = 0.01845510530133335

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp7"
$ws.Range("C4").Value = "Bmpr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03443933333333333
$ws.Range("H4").Value = 0.103318
$ws.Range("I4").Value = 0.05823261822459219
$ws.Range("J4").Value = 0.0582326182245922
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 38.53544233333333
$ws.Range("N4").Value = 115.606327
$ws.Range("O4").Value = 0.2608640200510233
$ws.Range("P4").Value = 0.2608640200510233
$ws.Range("Q4").Value = 1.327134943665111
$ws.Range("R4").Value = 11.944214492986
$ws.Range("S4").Value = 0.0151907948881636
$ws.Range("T4").Value = 0.01519079488816361

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Bmp7"
$ws.Range("C5").Value = "Bmpr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03443933333333333
$ws.Range("H5").Value = 0.103318
$ws.Range("I5").Value = 0.05823261822459219
$ws.Range("J5").Value = 0.0582326182245922
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 16.07945366666667
$ws.Range("N5").Value = 48.238361
$ws.Range("O5").Value = 0.1088491702632547
$ws.Range("P5").Value = 0.1088491702632547
$ws.Range("Q5").Value = 0.5537656646442222
$ws.Range("R5").Value = 4.983890981797999
$ws.Range("S5").Value = 0.006338572176003746
$ws.Range("T5").Value = 0.006338572176003747

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp7"
$ws.Range("C6").Value = "Bmpr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.51625
$ws.Range("H6").Value = 1.54875
$ws.Range("I6").Value = 0.8729143757654733
$ws.Range("J6").Value = 0.8729143757654734
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 46.29121633333333
$ws.Range("N6").Value = 138.873649
$ws.Range("O6").Value = 0.3133663986859022
$ws.Range("P6").Value = 0.3133663986859022
$ws.Range("Q6").Value = 23.89784043208333
$ws.Range("R6").Value = 215.08056388875
$ws.Range("S6").Value = 0.2735420342947787
$ws.Range("T6").Value = 0.2735420342947788

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp7"
$ws.Range("C7").Value = "Bmpr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.51625
$ws.Range("H7").Value = 1.54875
$ws.Range("I7").Value = 0.8729143757654733
$ws.Range("J7").Value = 0.8729143757654734
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 46.81622333333333
$ws.Range("N7").Value = 140.44867
$ws.Range("O7").Value = 0.3169204109998198
$ws.Range("P7").Value = 0.3169204109998198
$ws.Range("Q7").Value = 24.16887529583333
$ws.Range("R7").Value = 217.5198776625
$ws.Range("S7").Value = 0.2766443827352449
$ws.Range("T7").Value = 0.2766443827352449

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Bmp7"
$ws.Range("C8").Value = "Bmpr2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.51625
$ws.Range("H8").Value = 1.54875
$ws.Range("I8").Value = 0.8729143757654733
$ws.Range("J8").Value = 0.8729143757654734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 38.53544233333333
$ws.Range("N8").Value = 115.606327
$ws.Range("O8").Value = 0.2608640200510233
$ws.Range("P8").Value = 0.2608640200510233
$ws.Range("Q8").Value = 19.89392210458333
$ws.Range("R8").Value = 179.04529894125
$ws.Range("S8").Value = 0.2277119532225109
$ws.Range("T8").Value = 0.2277119532225109

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Bmp7"
$ws.Range("C9").Value = "Bmpr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.51625
$ws.Range("H9").Value = 1.54875
$ws.Range("I9").Value = 0.8729143757654733
$ws.Range("J9").Value = 0.8729143757654734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 16.07945366666667
$ws.Range("N9").Value = 48.238361
$ws.Range("O9").Value = 0.1088491702632547
$ws.Range("P9").Value = 0.1088491702632547
$ws.Range("Q9").Value = 8.301017955416667
$ws.Range("R9").Value = 74.70916159875
$ws.Range("S9").Value = 0.09501600551293872
$ws.Range("T9").Value = 0.09501600551293873

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Bmp7"
$ws.Range("C10").Value = "Bmpr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.04072033333333334
$ws.Range("H10").Value = 0.122161
$ws.Range("I10").Value = 0.06885300600993445
$ws.Range("J10").Value = 0.06885300600993445
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 46.29121633333333
$ws.Range("N10").Value = 138.873649
$ws.Range("O10").Value = 0.3133663986859022
$ws.Range("P10").Value = 0.3133663986859022
$ws.Range("Q10").Value = 1.884993759498778
$ws.Range("R10").Value = 16.964943835489
$ws.Range("S10").Value = 0.02157621853203194
$ws.Range("T10").Value = 0.02157621853203194

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Bmp7"
$ws.Range("C11").Value = "Bmpr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.04072033333333334
$ws.Range("H11").Value = 0.122161
$ws.Range("I11").Value = 0.06885300600993445
$ws.Range("J11").Value = 0.06885300600993445
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 46.81622333333333
$ws.Range("N11").Value = 140.44867
$ws.Range("O11").Value = 0.3169204109998198
$ws.Range("P11").Value = 0.3169204109998198
$ws.Range("Q11").Value = 1.906372219541111
$ws.Range("R11").Value = 17.15734997587
$ws.Range("S11").Value = 0.02182092296324149
$ws.Range("T11").Value = 0.02182092296324149

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Bmp7"
$ws.Range("C12").Value = "Bmpr2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.04072033333333334
$ws.Range("H12").Value = 0.122161
$ws.Range("I12").Value = 0.06885300600993445
$ws.Range("J12").Value = 0.06885300600993445
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 38.53544233333333
$ws.Range("N12").Value = 115.606327
$ws.Range("O12").Value = 0.2608640200510233
$ws.Range("P12").Value = 0.2608640200510233
$ws.Range("Q12").Value = 1.569176056960778
$ws.Range("R12").Value = 14.122584512647
$ws.Range("S12").Value = 0.01796127194034877
$ws.Range("T12").Value = 0.01796127194034877

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Bmp7"
$ws.Range("C13").Value = "Bmpr2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.04072033333333334
$ws.Range("H13").Value = 0.122161
$ws.Range("I13").Value = 0.06885300600993445
$ws.Range("J13").Value = 0.06885300600993445
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 16.07945366666667
$ws.Range("N13").Value = 48.238361
$ws.Range("O13").Value = 0.1088491702632547
$ws.Range("P13").Value = 0.1088491702632547
$ws.Range("Q13").Value = 0.6547607131245556
$ws.Range("R13").Value = 5.892846418121
$ws.Range("S13").Value = 0.007494592574312256
$ws.Range("T13").Value = 0.007494592574312256

Write-Host "Done writing Bmp7-Bmpr2 TPM update"
